$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Title text boxes: "PRO-FORMA INVOICE" -> "TAX INVOICE"
#    (lives twice in the same paragraph: DrawingML Choice + VML Fallback)
# ---------------------------------------------------------------
$pTitle = $d.Paragraphs.Item(5)
$xmlTitle = $pTitle.Range.WordOpenXML
$xmlTitle = $xmlTitle.Replace("PRO-FORMA INVOICE", "TAX INVOICE")
$pTitle.Range.InsertXML($xmlTitle)

# ---------------------------------------------------------------
# 2. "And choose carry." -> "Consider recognize."
#    (also duplicated Choice/Fallback content in one paragraph)
# ---------------------------------------------------------------
$pDesc = $d.Paragraphs.Item(24)
$xmlDesc = $pDesc.Range.WordOpenXML
$xmlDesc = $xmlDesc.Replace("And choose carry.", "Consider recognize.")
$pDesc.Range.InsertXML($xmlDesc)

# ---------------------------------------------------------------
# 3. Invoice header table (Table 1): Invoice #, Date, Customer #
# ---------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 1).Range.Text = "8"
$t1.Cell(2, 2).Range.Text = "1996-04-20"
$t1.Cell(4, 1).Range.Text = "7"

# ---------------------------------------------------------------
# 4. Bill-to block (plain paragraphs, single runs)
# ---------------------------------------------------------------
$d.Content.Find.Execute("Joseph LLC", $false, $false, $false, $false, $false, $true, 1, $false, "Charles Floyd", 2)
$d.Content.Find.Execute("07873 Harper Road Suite 134", $false, $false, $false, $false, $false, $true, 1, $false, "44442 Dawn Corners Apt. 118", 2)
$d.Content.Find.Execute("Williamschester, WI 20543", $false, $false, $false, $false, $false, $true, 1, $false, "Joneschester, VI 73504", 2)
$d.Content.Find.Execute("Wongshire", $false, $false, $false, $false, $false, $true, 1, $false, "West Lisa", 2)
$d.Content.Find.Execute("Luxembourg", $false, $false, $false, $false, $false, $true, 1, $false, "Paraguay", 2)

# ---------------------------------------------------------------
# 5. Line-items table (Table 2)
# ---------------------------------------------------------------
$t2 = $d.Tables.Item(2)

# Row 2 ("media")
$t2.Cell(2, 1).Range.Text = "progress"
$t2.Cell(2, 2).Range.Text = "White bring her well a century door too method language alone than."
$t2.Cell(2, 3).Range.Text = "28"
$t2.Cell(2, 4).Range.Find.Execute("494.52", $false, $false, $false, $false, $false, $true, 1, $false, "945.22", 2)
$t2.Cell(2, 5).Range.Find.Execute("29176.68", $false, $false, $false, $false, $false, $true, 1, $false, "26466.16", 2)

# Row 3 ("soil")
$t2.Cell(3, 1).Range.Text = "shock"
$t2.Cell(3, 2).Range.Text = "Plan sense far world either doctor statement."
$t2.Cell(3, 3).Range.Text = "2"
$t2.Cell(3, 4).Range.Find.Execute("71.74", $false, $false, $false, $false, $false, $true, 1, $false, "5548.94", 2)
$t2.Cell(3, 5).Range.Find.Execute("6097.90", $false, $false, $false, $false, $false, $true, 1, $false, "11097.88", 2)

# Row 4 ("stable")
$t2.Cell(4, 1).Range.Text = "photo"
$t2.Cell(4, 2).Range.Text = "Property court case recently activity understand so blood team find name work."
$t2.Cell(4, 3).Range.Text = "66"
$t2.Cell(4, 4).Range.Find.Execute("2043.27", $false, $false, $false, $false, $false, $true, 1, $false, "660.24", 2)
$t2.Cell(4, 5).Range.Find.Execute("89903.88", $false, $false, $false, $false, $false, $true, 1, $false, "43575.84", 2)

# Totals (still at original row numbers, before the row-5 delete below)
$t2.Cell(6, 5).Range.Text = "N$81,139.88"
$t2.Cell(7, 5).Range.Find.Execute("16.51", $false, $false, $false, $false, $false, $true, 1, $false, "15.0", 2)
$t2.Cell(8, 5).Range.Text = "N$1,217,098.20"
$t2.Cell(9, 5).Range.Text = "N$1,298,238.08"

# Row 5 ("studio") is removed entirely
$t2.Rows.Item(5).Delete()
